# Updated cryptos list (prices & volume %) as published by the
# "Updated cryptos list ... with GitHub Actions" commit.
# Note: price cells in column D are plain decimal-looking text
# (e.g. "357.23"), so NumberFormat is forced to "@" (Text) before
# assigning those values to stop Excel from auto-converting them to
# floating point numbers and corrupting the exact displayed string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '51.950.03'
$ws.Cells.Item(2, 5).Value = '  +0.41%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.819.54'
$ws.Cells.Item(3, 5).Value = '  +1.50%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.02%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '357.23'
$ws.Cells.Item(5, 5).Value = '  +0.10%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '110.34'
$ws.Cells.Item(6, 5).Value = '  +1.17%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.559'
$ws.Cells.Item(7, 5).Value = '  +0.59%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  -0.01%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.637'
$ws.Cells.Item(9, 5).Value = '  +8.55%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '40.28'
$ws.Cells.Item(10, 5).Value = '  +1.28%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  +0.27%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.0843'
$ws.Cells.Item(12, 5).Value = '  -0.23%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '20.08'
$ws.Cells.Item(13, 5).Value = '  +3.00%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '7.85'
$ws.Cells.Item(14, 5).Value = '  +2.99%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '3.257.35'
$ws.Cells.Item(15, 5).Value = '  +1.46%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '2.829.02'
$ws.Cells.Item(16, 5).Value = '  +1.58%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.947'
$ws.Cells.Item(17, 5).Value = '  +1.40%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '51.920.99'

# Row 19
$ws.Cells.Item(19, 5).Value = '  +3.07%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '3.19'
$ws.Cells.Item(20, 5).Value = '  +3.65%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '13.71'
$ws.Cells.Item(21, 5).Value = '  +4.16%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '0.0₃0981'
$ws.Cells.Item(22, 5).Value = '  +1.20%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '70.53'
$ws.Cells.Item(23, 5).Value = '  +0.47%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '269.04'
$ws.Cells.Item(24, 5).Value = '  +0.21%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  +1.21%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '26.26'
$ws.Cells.Item(26, 5).Value = '  -0.50%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  +0.02%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '0.165'
$ws.Cells.Item(28, 5).Value = '  +0.61%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '10.43'
$ws.Cells.Item(29, 5).Value = '  +1.98%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '38.26'
$ws.Cells.Item(30, 5).Value = '  +9.26%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  +0.77%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '6.22'
$ws.Cells.Item(32, 5).Value = '  -0.93%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '52.17'
$ws.Cells.Item(33, 5).Value = '  +0.72%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '5.70'
$ws.Cells.Item(34, 5).Value = '  +11.03%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.0449'
$ws.Cells.Item(35, 5).Value = '  -0.18%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.0871'
$ws.Cells.Item(36, 5).Value = '  +3.78%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.999'
$ws.Cells.Item(37, 5).Value = '  +0.01%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  +1.21%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '2.02'
$ws.Cells.Item(39, 5).Value = '  +2.93%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '3.16'
$ws.Cells.Item(40, 5).Value = '  +0.92%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  +1.26%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '2.52'
$ws.Cells.Item(42, 5).Value = '  -0.60%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '22.04'
$ws.Cells.Item(43, 5).Value = '  +1.41%  '

# Row 44
$ws.Cells.Item(44, 2).Value = 'WEMIXToken'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '2.20'
$ws.Cells.Item(44, 5).Value = '  -1.19%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'Monero'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '119.43'
$ws.Cells.Item(45, 5).Value = '  +0.02%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '2.49'
$ws.Cells.Item(46, 5).Value = '  +9.05%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'NEARProtocol'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '3.41'
$ws.Cells.Item(47, 5).Value = '  +3.92%  '

# Row 48
$ws.Cells.Item(48, 2).Value = 'Maker'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(48, 4).Value = '2.115.22'
$ws.Cells.Item(48, 5).Value = '  +1.39%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.934'
$ws.Cells.Item(49, 5).Value = '  -1.03%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '1.38'
$ws.Cells.Item(50, 5).Value = '  +9.72%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '5.47'
$ws.Cells.Item(51, 5).Value = '  -3.01%  '
